$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Responsable_etapa"

# 2. Column widths
$ws.Columns.Item(1).ColumnWidth = 10.75
$ws.Columns.Item(2).ColumnWidth = 12.75
$ws.Columns.Item(3).ColumnWidth = 6.75
$ws.Columns.Item(4).ColumnWidth = 12.75
$ws.Columns.Item(5).ColumnWidth = 14.75
$ws.Columns.Item(6).ColumnWidth = 16.75
$ws.Columns.Item(7).ColumnWidth = 18.75
$ws.Columns.Item(8).ColumnWidth = 8.75
$ws.Columns.Item(9).ColumnWidth = 9.75
$ws.Columns.Item(10).ColumnWidth = 11.75
$ws.Columns.Item(11).ColumnWidth = 21.75
$ws.Columns.Item(12).ColumnWidth = 7.75
$ws.Columns.Item(13).ColumnWidth = 7.75
$ws.Columns.Item(14).ColumnWidth = 7.75
$ws.Columns.Item(15).ColumnWidth = 7.75
$ws.Columns.Item(16).ColumnWidth = 7.75
$ws.Columns.Item(17).ColumnWidth = 7.75

# 3. View options: hide gridlines + freeze header row
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 4. Header row formatting (row 1): bold white font, dark blue fill, bottom border, centered
$header = $ws.Range("A1:Q1")
$header.Font.Color = 16777215
$header.Font.Bold = $true
$header.Interior.Color = 7949855
$header.Borders.Item(9).Color = 0
$header.Borders.Item(9).LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108

# 5. Data rows formatting (rows 2-20)
$dataAll = $ws.Range("A2:Q20")
$dataAll.Font.Name = "Calibri"
$dataAll.Borders.Color = 0
$dataAll.Borders.LineStyle = 1
$dataAll.VerticalAlignment = -4108

# 5a. Text columns (A,B) -> left aligned
$ws.Range("A2:B20").HorizontalAlignment = -4131

# 5b. Center-aligned general columns (L:Q)
$ws.Range("L2:Q20").HorizontalAlignment = -4108

# 5c. Integer columns (C,D,F,H,I,K) -> center aligned + #,##0 format
$intCols = "C2:C20,D2:D20,F2:F20,H2:H20,I2:I20,K2:K20"
foreach ($a in $intCols.Split(",")) {
  $r = $ws.Range($a)
  $r.HorizontalAlignment = -4108
  $r.NumberFormat = "#,##0"
}

# 5d. Percentage columns (E,G,J) -> center aligned + 0.0"%" format
$pctCols = "E2:E20,G2:G20,J2:J20"
foreach ($a in $pctCols.Split(",")) {
  $r = $ws.Range($a)
  $r.HorizontalAlignment = -4108
  $r.NumberFormat = "0.0""%"""
}

# 6. AutoFilter over the full data range
$ws.Range("A1:Q20").AutoFilter()

# 7. Hidden defined name backing the filter (mirrors Excel's _FilterDatabase)
$fdbName = $ws.Names.Add("_xlnm._FilterDatabase", "='Responsable_etapa'!`$A`$1:`$Q`$20")
$fdbName.Visible = $false

Write-Host "step done"
